$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.878.86"
$ws.Range("E2").Value = "  +6.06%  "

$ws.Range("D3").Value = "3.510.28"
$ws.Range("E3").Value = "  +6.74%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "191.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +10.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "557.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.75%  "

$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D7").Value = "3.512.25"
$ws.Range("E7").Value = "  +7.02%  "

$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.615"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.49%  "

$ws.Range("E9").Value = "  -0.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.642"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.50%  "

$ws.Range("E12").Value = "  +13.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000277"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.69%  "

$ws.Range("D15").Value = "4.052.09"
$ws.Range("E15").Value = "  +6.13%  "

$ws.Range("D16").Value = "3.497.98"
$ws.Range("E16").Value = "  +6.31%  "

$ws.Range("D17").Value = "68.245.04"
$ws.Range("E17").Value = "  +6.73%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.51%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.121"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.78%  "

$ws.Range("E21").Value = "  +6.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "406.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.30%  "

$ws.Range("E26").Value = "  +8.08%  "

$ws.Range("E27").Value = "  +10.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.66%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "688.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.21%  "

$ws.Range("E35").Value = "  +6.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "60.81"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.40%  "

$ws.Range("D37").Value = "0.0₃0837"
$ws.Range("E37").Value = "  +21.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "39.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.406"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.30%  "

$ws.Range("E40").Value = "  +0.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +25.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +16.48%  "

$ws.Range("E43").Value = "  +11.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.17%  "

$ws.Range("D45").Value = "3.052.62"
$ws.Range("E45").Value = "  +4.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0425"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +14.10%  "

$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +13.05%  "

